$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - rename existing headers and add two new ones
$ws.Range("A1").Value = "اسم المادة"
$ws.Range("B1").Value = "كود المادة"
$ws.Range("C1").Value = "القسم"
$ws.Range("D1").Value = "سنة دراسية"

# Give the two new header cells the same bold styling as the rest of the header row
$ws.Range("C1").Font.Bold = $true
$ws.Range("D1").Font.Bold = $true

# Row 2 - first course record
$ws.Range("A2").Value = "Electronics"
$ws.Range("B2").Value = "ECE213"
$ws.Range("C2").Value = "عام"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1"
$ws.Range("D2").ClearFormats()

# Row 3 - second course record (replaces old row 3 + row 4 data)
$ws.Range("A3").Value = "Electronics"
$ws.Range("B3").Value = "ECE213"
$ws.Range("C3").Value = "اتصالات"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1"
$ws.Range("D3").ClearFormats()

# The old data occupied a 4th row; clear it since the new data only spans 3 rows
$ws.Range("A4:D4").ClearContents()
